# Fruta / hortaliza, semanal
# Insert 3 new rows before row 706 (pushing the existing data for rows
# 706-726 down to 709-729) and populate the new rows with the "Early John"
# Nectarin entries for Provincia de San Felipe de Aconcagua.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows above the current row 706, shifting everything
# below (old rows 706-726) down to rows 709-729.
$ws.Rows.Item(706).Insert()
$ws.Rows.Item(706).Insert()
$ws.Rows.Item(706).Insert()

# New row 706
$ws.Cells.Item(706, 1).Value = 9
$ws.Cells.Item(706, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(706, 3).Value = "Metropolitana"
$ws.Cells.Item(706, 4).Value = 44890
$ws.Cells.Item(706, 5).Value = 13
$ws.Cells.Item(706, 6).Value = "Fruta"
$ws.Cells.Item(706, 7).Value = 100103
$ws.Cells.Item(706, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(706, 9).Value = 100103006
$ws.Cells.Item(706, 10).Value = "Nectarín"
$ws.Cells.Item(706, 11).Value = "Early John"
$ws.Cells.Item(706, 12).Value = "Especial"
$ws.Cells.Item(706, 13).Value = 170
$ws.Cells.Item(706, 14).Value = 10000
$ws.Cells.Item(706, 15).Value = 10000
$ws.Cells.Item(706, 16).Value = 10000
$ws.Cells.Item(706, 17).Value = "$/bandeja 8 kilos empedrada"
$ws.Cells.Item(706, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(706, 19).Value = 1250
$ws.Cells.Item(706, 20).Value = 8

# New row 707
$ws.Cells.Item(707, 1).Value = 9
$ws.Cells.Item(707, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(707, 3).Value = "Metropolitana"
$ws.Cells.Item(707, 4).Value = 44890
$ws.Cells.Item(707, 5).Value = 13
$ws.Cells.Item(707, 6).Value = "Fruta"
$ws.Cells.Item(707, 7).Value = 100103
$ws.Cells.Item(707, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(707, 9).Value = 100103006
$ws.Cells.Item(707, 10).Value = "Nectarín"
$ws.Cells.Item(707, 11).Value = "Early John"
$ws.Cells.Item(707, 12).Value = "Primera"
$ws.Cells.Item(707, 13).Value = 300
$ws.Cells.Item(707, 14).Value = 8000
$ws.Cells.Item(707, 15).Value = 8000
$ws.Cells.Item(707, 16).Value = 8000
$ws.Cells.Item(707, 17).Value = "$/bandeja 8 kilos empedrada"
$ws.Cells.Item(707, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(707, 19).Value = 1000
$ws.Cells.Item(707, 20).Value = 8

# New row 708
$ws.Cells.Item(708, 1).Value = 9
$ws.Cells.Item(708, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(708, 3).Value = "Metropolitana"
$ws.Cells.Item(708, 4).Value = 44890
$ws.Cells.Item(708, 5).Value = 13
$ws.Cells.Item(708, 6).Value = "Fruta"
$ws.Cells.Item(708, 7).Value = 100103
$ws.Cells.Item(708, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(708, 9).Value = 100103006
$ws.Cells.Item(708, 10).Value = "Nectarín"
$ws.Cells.Item(708, 11).Value = "Early John"
$ws.Cells.Item(708, 12).Value = "Segunda"
$ws.Cells.Item(708, 13).Value = 280
$ws.Cells.Item(708, 14).Value = 6000
$ws.Cells.Item(708, 15).Value = 6000
$ws.Cells.Item(708, 16).Value = 6000
$ws.Cells.Item(708, 17).Value = "$/bandeja 8 kilos empedrada"
$ws.Cells.Item(708, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(708, 19).Value = 750
$ws.Cells.Item(708, 20).Value = 8

# Re-apply the date number format (style index 2 in the original workbook)
# to the D column of the newly inserted rows, matching the rest of the
# column (Excel's row-insert already inherits this from row 706/709, but
# set it explicitly to be safe).
$ws.Range("D706:D708").NumberFormat = "YYYY-MM-DD HH:MM:SS"
